# Auto-generated edit script applying the Phantom_Profits diff
# Updates numeric cells across ALC, ARM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3013.52
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 3013.52
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 9040.559999999999
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -9376.559999999999
$ws.Range("H21").Value = 26249.75
$ws.Range("I21").Value = 26249.75
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 26249.75
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -25781.75
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 26249.75
$ws.Range("I23").Value = 26249.75
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 26249.75
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -26015.75
$ws.Range("N23").ClearContents()
$ws.Range("H28").Value = 1344.875
$ws.Range("I28").Value = 1437.4546
$ws.Range("J28").Value = 1141.2
$ws.Range("K28").Value = 1437.4546
$ws.Range("L28").Value = 1141.2
$ws.Range("M28").Value = -952.4546
$ws.Range("N28").Value = -2111.2
$ws.Range("H40").Value = 1242.5834
$ws.Range("I40").Value = 1242.5834
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1242.5834
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1067.5834
$ws.Range("N40").ClearContents()
$ws.Range("H62").Value = 1158.2
$ws.Range("I62").Value = 1158.2
$ws.Range("K62").Value = 1158.2
$ws.Range("M62").Value = -534.2
$ws.Range("H65").Value = 1158.2
$ws.Range("I65").Value = 1158.2
$ws.Range("K65").Value = 5791
$ws.Range("M65").Value = -2671
$ws.Range("H69").Value = 10665
$ws.Range("J69").Value = 9998.5
$ws.Range("L69").Value = 29995.5
$ws.Range("N69").Value = -31743.5
$ws.Range("H72").Value = 10665
$ws.Range("J72").Value = 9998.5
$ws.Range("L72").Value = 89986.5
$ws.Range("N72").Value = -98722.5
$ws.Range("H86").Value = 9761.200000000001
$ws.Range("I86").Value = 14999
$ws.Range("J86").Value = 8451.75
$ws.Range("K86").Value = 14999
$ws.Range("L86").Value = 8451.75
$ws.Range("M86").Value = -13876
$ws.Range("N86").Value = -10697.75
$ws.Range("H89").Value = 9761.200000000001
$ws.Range("I89").Value = 14999
$ws.Range("J89").Value = 8451.75
$ws.Range("K89").Value = 74995
$ws.Range("L89").Value = 42258.75
$ws.Range("M89").Value = -69379
$ws.Range("N89").Value = -53490.75
$ws.Range("H106").Value = 9197.799999999999
$ws.Range("I106").Value = 10748
$ws.Range("J106").Value = 2997
$ws.Range("K106").Value = 10748
$ws.Range("L106").Value = 2997
$ws.Range("M106").Value = -10117
$ws.Range("N106").Value = -4259
$ws.Range("H132").Value = 4453.3887
$ws.Range("I132").Value = 4538.8823
$ws.Range("K132").Value = 13616.6469
$ws.Range("M132").Value = -11086.6469
$ws.Range("H135").Value = 1029.1818
$ws.Range("I135").Value = 1100.1
$ws.Range("K135").Value = 9900.9
$ws.Range("M135").Value = -7365.9
$ws.Range("H137").Value = 4651.0835
$ws.Range("I137").Value = 3961.3
$ws.Range("K137").Value = 11883.9
$ws.Range("M137").Value = -9333.900000000001
$ws.Range("H141").Value = 2221.5334
$ws.Range("I141").Value = 2221.5334
$ws.Range("K141").Value = 6664.600199999999
$ws.Range("M141").Value = -1484.600199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2867.3333
$ws.Range("I45").Value = 3240.8
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 3240.8
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -2863.8
$ws.Range("N45").Value = -1754
$ws.Range("H61").Value = 4164.4546
$ws.Range("I61").Value = 3780.9
$ws.Range("K61").Value = 3780.9
$ws.Range("M61").Value = -3568.9
$ws.Range("H74").Value = 2260
$ws.Range("I74").Value = 2140.625
$ws.Range("K74").Value = 2140.625
$ws.Range("M74").Value = -1266.625
$ws.Range("H77").Value = 2260
$ws.Range("I77").Value = 2140.625
$ws.Range("K77").Value = 10703.125
$ws.Range("M77").Value = -6335.125
$ws.Range("H132").Value = 3960.4
$ws.Range("I132").Value = 3960.4
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11881.2
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9351.200000000001
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 4164.4546
$ws.Range("I136").Value = 3780.9
$ws.Range("K136").Value = 11342.7
$ws.Range("M136").Value = -8792.700000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4123
$ws.Range("I31").Value = 2066.6667
$ws.Range("K31").Value = 2066.6667
$ws.Range("M31").Value = -1771.6667
$ws.Range("H34").Value = 4123
$ws.Range("I34").Value = 2066.6667
$ws.Range("K34").Value = 2066.6667
$ws.Range("M34").Value = -1864.6667
$ws.Range("H37").Value = 24991.4
$ws.Range("I37").Value = 24994.5
$ws.Range("J37").Value = 24989.334
$ws.Range("K37").Value = 24994.5
$ws.Range("L37").Value = 24989.334
$ws.Range("M37").Value = -24887.5
$ws.Range("N37").Value = -25203.334
$ws.Range("H58").Value = 5966.6665
$ws.Range("I58").Value = 3950
$ws.Range("K58").Value = 3950
$ws.Range("M58").Value = -3747
$ws.Range("H63").Value = 72297.82000000001
$ws.Range("J63").Value = 81697.336
$ws.Range("L63").Value = 81697.336
$ws.Range("N63").Value = -83069.336
$ws.Range("H66").Value = 72297.82000000001
$ws.Range("J66").Value = 81697.336
$ws.Range("L66").Value = 245092.008
$ws.Range("N66").Value = -251956.008
$ws.Range("H86").Value = 8999.666999999999
$ws.Range("I86").Value = 8999.5
$ws.Range("K86").Value = 8999.5
$ws.Range("M86").Value = -7876.5
$ws.Range("H89").Value = 8999.666999999999
$ws.Range("I89").Value = 8999.5
$ws.Range("K89").Value = 44997.5
$ws.Range("M89").Value = -39381.5
$ws.Range("H99").Value = 2995.25
$ws.Range("I99").Value = 2995.25
$ws.Range("K99").Value = 2995.25
$ws.Range("M99").Value = -1497.25
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H122").Value = 7168.6665
$ws.Range("I122").Value = 8999.5
$ws.Range("J122").Value = 3507
$ws.Range("K122").Value = 26998.5
$ws.Range("L122").Value = 10521
$ws.Range("M122").Value = -24548.5
$ws.Range("N122").Value = -15421
$ws.Range("H126").Value = 2995.25
$ws.Range("I126").Value = 2995.25
$ws.Range("K126").Value = 8985.75
$ws.Range("M126").Value = -6515.75
$ws.Range("H134").Value = 4524.3
$ws.Range("I134").Value = 4471.4443
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 13414.3329
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -10879.3329
$ws.Range("N134").Value = -20070
$ws.Range("H136").Value = 5966.6665
$ws.Range("I136").Value = 3950
$ws.Range("K136").Value = 11850
$ws.Range("M136").Value = -9300

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H36").Value = 5002
$ws.Range("I36").Value = 5002
$ws.Range("K36").Value = 15006
$ws.Range("M36").Value = -14837
$ws.Range("H51").Value = 799
$ws.Range("I51").Value = 799
$ws.Range("K51").Value = 2397
$ws.Range("M51").Value = -1937
$ws.Range("H131").Value = 2049.875
$ws.Range("I131").Value = 2133
$ws.Range("J131").Value = 2000
$ws.Range("K131").Value = 6399
$ws.Range("L131").Value = 6000
$ws.Range("M131").Value = -1359
$ws.Range("N131").Value = -16080

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2134.8333
$ws.Range("I80").Value = 2040.6
$ws.Range("J80").Value = 2606
$ws.Range("K80").Value = 2040.6
$ws.Range("L80").Value = 2606
$ws.Range("M80").Value = -1042.6
$ws.Range("N80").Value = -4602
$ws.Range("H83").Value = 2134.8333
$ws.Range("I83").Value = 2040.6
$ws.Range("J83").Value = 2606
$ws.Range("K83").Value = 10203
$ws.Range("L83").Value = 13030
$ws.Range("M83").Value = -5211
$ws.Range("N83").Value = -23014

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1719.7142
$ws.Range("J46").Value = 1835
$ws.Range("L46").Value = 1835
$ws.Range("N46").Value = -2211
$ws.Range("H82").Value = 2069.6365
$ws.Range("I82").Value = 1487.5
$ws.Range("J82").Value = 2402.2856
$ws.Range("K82").Value = 1487.5
$ws.Range("L82").Value = 2402.2856
$ws.Range("M82").Value = -1126.5
$ws.Range("N82").Value = -3124.2856
$ws.Range("H85").Value = 2069.6365
$ws.Range("I85").Value = 1487.5
$ws.Range("J85").Value = 2402.2856
$ws.Range("K85").Value = 1487.5
$ws.Range("L85").Value = 2402.2856
$ws.Range("M85").Value = -239.5
$ws.Range("N85").Value = -4898.2856

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4230.5557
$ws.Range("I126").Value = 3884.375
$ws.Range("K126").Value = 11653.125
$ws.Range("M126").Value = -9183.125
